# v0.3.1 release edits: focus dashboards / amount buckets / popover UX
# Applies the row 50/51 follow-up fields plus the brand-new row 52
# (CPOE / AMISTAD project) on the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50 (CCDC-86 / 秘鲁58区 follow-up): fill in the previously blank
#     project-name / nickname / contract-number / remark cells.
$ws.Range("J50").Value = "秘鲁58区"
$ws.Range("K50").Value = "58区修井"
$ws.Range("L50").Value = "8"
$ws.Range("V50").Value = "长钻50090"

# --- Row 51 (CCDC-87): no textual change (its shared strings simply shift
#     index because of the inserts above/below); nothing to write here.

# --- Row 52 (new): CPOE / 中油海工 / AMISTAD project entry.
$ws.Range("A52").Value = "厄瓜多尔"
$ws.Range("B52").Value = 15
$ws.Range("C52").Value = "CPOE"
$ws.Range("I52").Value = "中油海工"
$ws.Range("J52").Value = "AMISTAD项目"
$ws.Range("K52").Value = "海上修井"
$ws.Range("L52").Value = "9"
$ws.Range("Q52").Value = "动搬迁"
$ws.Range("R52").Value = "目前：整改设备，作业前准备，配合甲方测试设备，配钻具。"
$ws.Range("V52").Value = "CPOE"

# --- Sheet view: scroll/selection moved down to the newly-edited area.
$ws.Range("L50").Select()
